# Update "想去人数" (want-to-go count) values in column F across the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets to
# match the refreshed data snapshot published to gh-pages.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetShow    = $wb.Worksheets.Item("演出")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# Row => new value updates for the "展览" sheet (sheet1)
$exhibitUpdates = @{
    5  = 78
    7  = 596
    8  = 116
    9  = 8777
    11 = 330
    13 = 993
    18 = 267
    19 = 68
    20 = 233
    21 = 1050
}

foreach ($row in $exhibitUpdates.Keys) {
    $sheetExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Row => new value updates for the "演出" sheet (sheet2)
$showUpdates = @{
    2 = 2
}

foreach ($row in $showUpdates.Keys) {
    $sheetShow.Cells.Item($row, 6).Value = $showUpdates[$row]
}

# Row => new value updates for the "全部类型" sheet (sheet4)
$allUpdates = @{
    5  = 2
    6  = 78
    9  = 596
    10 = 116
    11 = 8777
    13 = 330
    15 = 993
    20 = 267
    21 = 68
    22 = 233
    23 = 1050
}

foreach ($row in $allUpdates.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
